$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Row 13 - Quadrato (square) of pairwise sums, mirroring columns F:K
$ws.Range("F13").Formula = "=(B1+B2)^2"
$ws.Range("G13").Formula = "=(B1+B3)^2"
$ws.Range("H13").Formula = "=(B1+B4)^2"
$ws.Range("I13").Formula = "=(B2+B3)^2"
$ws.Range("J13").Formula = "=(B2+B4)^2"
$ws.Range("K13").Formula = "=(B3+B4)^2"

# Row 14 - Cubo (cube) of pairwise sums, mirroring columns F:K
$ws.Range("F14").Formula = "=(B1+B2)^3"
$ws.Range("G14").Formula = "=(B1+B3)^3"
$ws.Range("H14").Formula = "=(B1+B4)^3"
$ws.Range("I14").Formula = "=(B2+B3)^3"
$ws.Range("J14").Formula = "=(B2+B4)^3"
$ws.Range("K14").Formula = "=(B3+B4)^3"

# Row 15 - Radice Quadrata (square root) of pairwise sums, mirroring columns F:K
$ws.Range("F15").Formula = "=SQRT(B1+B2)"
$ws.Range("G15").Formula = "=SQRT(B1+B3)"
$ws.Range("H15").Formula = "=SQRT(B1+B4)"
$ws.Range("I15").Formula = "=SQRT(B2+B3)"
$ws.Range("J15").Formula = "=SQRT(B2+B4)"
$ws.Range("K15").Formula = "=SQRT(B3+B4)"

# Update the saved selection to K15 as in the source workbook
$ws.Range("K15").Select()
